$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

$ws.Cells.Item($row, 1).Value = 42607.886458333334
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = -12
$ws.Cells.Item($row, 3).Value = 53
$ws.Cells.Item($row, 4).Value = 45
$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 83
$ws.Cells.Item($row, 7).Value = 17249
$ws.Cells.Item($row, 8).Value = 17426
$ws.Cells.Item($row, 9).Value = 986
$ws.Cells.Item($row, 10).Value = 174
$ws.Cells.Item($row, 11).Value = 146
$ws.Cells.Item($row, 12).Value = 3
$ws.Cells.Item($row, 13).Value = 15
$ws.Cells.Item($row, 14).Value = "Named"
